$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "sd"
$ws.Range("J2").Value = "Statement-non-opinion"
$ws.Range("I5").Value = "sd"
$ws.Range("J5").Value = "Statement-non-opinion"
$ws.Range("I9").Value = "sv"
$ws.Range("J9").Value = "Statement-opinion"
$ws.Range("I16").Value = "sv"
$ws.Range("J16").Value = "Statement-opinion"
$ws.Range("I22").Value = "ba"
$ws.Range("J22").Value = "Appreciation"
$ws.Range("I27").Value = "ba"
$ws.Range("J27").Value = "Appreciation"
$ws.Range("I32").Value = "sv"
$ws.Range("J32").Value = "Statement-opinion"
$ws.Range("I45").Value = "sv"
$ws.Range("J45").Value = "Statement-opinion"
$ws.Range("I48").Value = "aa"
$ws.Range("J48").Value = "Agree/Accept"
$ws.Range("I52").Value = "aa"
$ws.Range("J52").Value = "Agree/Accept"
$ws.Range("I62").Value = "aa"
$ws.Range("J62").Value = "Agree/Accept"
$ws.Range("I82").Value = "aa"
$ws.Range("J82").Value = "Agree/Accept"
$ws.Range("I85").Value = "b"
$ws.Range("J85").Value = "Acknowledge (Backchannel)"
$ws.Range("I100").Value = "sd"
$ws.Range("J100").Value = "Statement-non-opinion"
$ws.Range("I104").Value = "sv"
$ws.Range("J104").Value = "Statement-opinion"
$ws.Range("I114").Value = "sv"
$ws.Range("J114").Value = "Statement-opinion"
$ws.Range("I123").Value = "sd"
$ws.Range("J123").Value = "Statement-non-opinion"
$ws.Range("I133").Value = "b"
$ws.Range("J133").Value = "Acknowledge (Backchannel)"
$ws.Range("I137").Value = "sv"
$ws.Range("J137").Value = "Statement-opinion"
$ws.Range("I138").Value = "sd"
$ws.Range("J138").Value = "Statement-non-opinion"
$ws.Range("I142").Value = "sd"
$ws.Range("J142").Value = "Statement-non-opinion"
$ws.Range("I144").Value = "sd"
$ws.Range("J144").Value = "Statement-non-opinion"
$ws.Range("I146").Value = "sv"
$ws.Range("J146").Value = "Statement-opinion"
$ws.Range("I148").Value = "aa"
$ws.Range("J148").Value = "Agree/Accept"
$ws.Range("I154").Value = "sd"
$ws.Range("J154").Value = "Statement-non-opinion"
$ws.Range("I167").Value = "sv"
$ws.Range("J167").Value = "Statement-opinion"
$ws.Range("I184").Value = "aa"
$ws.Range("J184").Value = "Agree/Accept"
$ws.Range("I193").Value = "sv"
$ws.Range("J193").Value = "Statement-opinion"
$ws.Range("I199").Value = "aa"
$ws.Range("J199").Value = "Agree/Accept"
$ws.Range("I210").Value = "b"
$ws.Range("J210").Value = "Acknowledge (Backchannel)"
$ws.Range("I218").Value = "qy"
$ws.Range("J218").Value = "Yes-No-Question"
$ws.Range("I240").Value = "aa"
$ws.Range("J240").Value = "Agree/Accept"
$ws.Range("I241").Value = "sd"
$ws.Range("J241").Value = "Statement-non-opinion"
$ws.Range("I251").Value = "sd"
$ws.Range("J251").Value = "Statement-non-opinion"
$ws.Range("I252").Value = "ba"
$ws.Range("J252").Value = "Appreciation"
$ws.Range("I274").Value = "sv"
$ws.Range("J274").Value = "Statement-opinion"
$ws.Range("I282").Value = "b"
$ws.Range("J282").Value = "Acknowledge (Backchannel)"
$ws.Range("I310").Value = "sd"
$ws.Range("J310").Value = "Statement-non-opinion"
$ws.Range("I331").Value = "ba"
$ws.Range("J331").Value = "Appreciation"
$ws.Range("I348").Value = "aa"
$ws.Range("J348").Value = "Agree/Accept"
$ws.Range("I351").Value = "sv"
$ws.Range("J351").Value = "Statement-opinion"
$ws.Range("I353").Value = "ba"
$ws.Range("J353").Value = "Appreciation"
$ws.Range("I357").Value = "sv"
$ws.Range("J357").Value = "Statement-opinion"
$ws.Range("I362").Value = "aa"
$ws.Range("J362").Value = "Agree/Accept"
$ws.Range("I365").Value = "aa"
$ws.Range("J365").Value = "Agree/Accept"
$ws.Range("I366").Value = "sv"
$ws.Range("J366").Value = "Statement-opinion"
$ws.Range("I368").Value = "ba"
$ws.Range("J368").Value = "Appreciation"
$ws.Range("I372").Value = "aa"
$ws.Range("J372").Value = "Agree/Accept"
$ws.Range("I378").Value = "aa"
$ws.Range("J378").Value = "Agree/Accept"
$ws.Range("I381").Value = "b"
$ws.Range("J381").Value = "Acknowledge (Backchannel)"
$ws.Range("I390").Value = "b"
$ws.Range("J390").Value = "Acknowledge (Backchannel)"
$ws.Range("I398").Value = "ba"
$ws.Range("J398").Value = "Appreciation"
$ws.Range("I404").Value = "b"
$ws.Range("J404").Value = "Acknowledge (Backchannel)"
$ws.Range("I421").Value = "ba"
$ws.Range("J421").Value = "Appreciation"
$ws.Range("I429").Value = "sd"
$ws.Range("J429").Value = "Statement-non-opinion"
$ws.Range("I435").Value = "sd"
$ws.Range("J435").Value = "Statement-non-opinion"
$ws.Range("I446").Value = "ba"
$ws.Range("J446").Value = "Appreciation"
$ws.Range("I453").Value = "b"
$ws.Range("J453").Value = "Acknowledge (Backchannel)"
$ws.Range("I455").Value = "%"
$ws.Range("J455").Value = "Uninterpretable"
$ws.Range("I457").Value = "b"
$ws.Range("J457").Value = "Acknowledge (Backchannel)"
$ws.Range("I475").Value = "sd"
$ws.Range("J475").Value = "Statement-non-opinion"
